$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their original Text storage type so that
# numeric-looking strings (e.g. "206.42", "2.16") are not auto-converted
# into numbers by Excel when we assign the new values.
$targetCells = @('D2','E2','D3','E3','E4','D5','E5','D6','E7','D8','E8','D9','E9','E10','D11','E11','D12','E12','D13','E13','E14','E15','D16','E16','D17','E17','D18','E18','E19','D20','E20','E21','D22','E22','D23','E23','E24','D25','E25','E26','D27','E27','E28','E29','E30','E31','E32','D33','E33','E34','D36','E36','E37','E38','E39','E40','E41','D42','E42','D43','E43','B44','C44','D44','E44','B45','C45','D45','E45','E46','D47','E47','D48','E48','D49','E49','D50','E50','E51')
foreach ($addr in $targetCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '27.135.71'
$ws.Range('E2').Value = '  -2.07%  '
$ws.Range('D3').Value = '1.558.83'
$ws.Range('E3').Value = '  -2.16%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '206.42'
$ws.Range('E5').Value = '  -1.01%  '
$ws.Range('D6').Value = '0.490'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '22.16'
$ws.Range('E8').Value = '  -0.78%  '
$ws.Range('D9').Value = '0.247'
$ws.Range('E9').Value = '  -2.21%  '
$ws.Range('E10').Value = '  -0.21%  '
$ws.Range('D11').Value = '0.0862'
$ws.Range('E11').Value = '  -0.81%  '
$ws.Range('D12').Value = '1.780.37'
$ws.Range('E12').Value = '  -2.11%  '
$ws.Range('D13').Value = '1.564.40'
$ws.Range('E13').Value = '  -1.78%  '
$ws.Range('E14').Value = '  -2.59%  '
$ws.Range('E15').Value = '  -3.12%  '
$ws.Range('D16').Value = '62.87'
$ws.Range('E16').Value = '  -0.96%  '
$ws.Range('D17').Value = '27.141.09'
$ws.Range('E17').Value = '  -1.98%  '
$ws.Range('D18').Value = '215.11'
$ws.Range('E18').Value = '  -2.34%  '
$ws.Range('E19').Value = '  -1.59%  '
$ws.Range('D20').Value = '7.22'
$ws.Range('E20').Value = '  -1.77%  '
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('D22').Value = '4.11'
$ws.Range('E22').Value = '  -1.24%  '
$ws.Range('D23').Value = '9.33'
$ws.Range('E23').Value = '  -3.61%  '
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('D25').Value = '151.91'
$ws.Range('E25').Value = '  -1.40%  '
$ws.Range('E26').Value = '  -3.32%  '
$ws.Range('D27').Value = '14.91'
$ws.Range('E27').Value = '  -1.68%  '
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('E29').Value = '  -1.54%  '
$ws.Range('E30').Value = '  -1.98%  '
$ws.Range('E31').Value = '  -2.08%  '
$ws.Range('E32').Value = '  -2.03%  '
$ws.Range('D33').Value = '1.379.80'
$ws.Range('E33').Value = '  +0.43%  '
$ws.Range('E34').Value = '  -0.95%  '
$ws.Range('D36').Value = '0.943'
$ws.Range('E36').Value = '  -3.11%  '
$ws.Range('E37').Value = '  -1.78%  '
$ws.Range('E38').Value = '  -1.65%  '
$ws.Range('E39').Value = '  -2.21%  '
$ws.Range('E40').Value = '  -4.73%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').Value = '0.991'
$ws.Range('E42').Value = '  +2.47%  '
$ws.Range('D43').Value = '1.79'
$ws.Range('E43').Value = '  +3.54%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = '63.26'
$ws.Range('E44').Value = '  -1.95%  '
$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D45').Value = '2.16'
$ws.Range('E45').Value = '  -0.29%  '
$ws.Range('E46').Value = '  +0.14%  '
$ws.Range('D47').Value = '1.693.20'
$ws.Range('E47').Value = '  -2.07%  '
$ws.Range('D48').Value = '85.41'
$ws.Range('E48').Value = '  -1.89%  '
$ws.Range('D49').Value = '0.0₇0984'
$ws.Range('E49').Value = '  -2.84%  '
$ws.Range('D50').Value = '0.0491'
$ws.Range('E50').Value = '  -0.86%  '
$ws.Range('E51').Value = '  +0.09%  '
